# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
# Updates MSME indicator figures on the Colombia Summary sheet with more
# precise decimal values, while keeping them stored as text (matching the
# original workbook layout where these figures are plain text strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Enterprises density (per 1000 people) row
Set-TextValue "B11" "30.94"
Set-TextValue "C11" "1.24"
Set-TextValue "D11" "32.18"

# Enterprises (% of total) row
Set-TextValue "B13" "96.01"
Set-TextValue "C13" "3.86"
Set-TextValue "D13" "99.87"

# SME Associations source - density-like row (13.8 / 0.6 / 14.4 -> 13.8 / 0.62 / 14.42)
Set-TextValue "C32" "0.62"
Set-TextValue "D32" "14.42"
